$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell U1 - same visual style as the other header cells
# (bold, bordered, centered) but stored as literal text, not an auto-detected date.
$ws.Range("U1").NumberFormat = "@"
$ws.Range("U1").Font.Bold = $true
$ws.Range("U1").Borders.LineStyle = 1
$ws.Range("U1").HorizontalAlignment = -4108
$ws.Range("U1").VerticalAlignment = -4160
$ws.Range("U1").Value = "2025-06-23"

# Update totals for row 2
$ws.Range("S2").Value = 16
$ws.Range("T2").Value = 6.2
$ws.Range("U2").Value = "❌"

# Update totals for row 3
$ws.Range("S3").Value = 16
$ws.Range("T3").Value = 6.2
$ws.Range("U3").Value = "❌"
